{"js": "// Move the \"_GoBack\" bookmark from the empty list item (after the\n// \"Charts we wanted to build using D3.js?\" line) down onto the empty\n// list item left behind after merging the two \"Heat map\" bullets\n// (\"... (by state)\" and \"Heat map for GA (by counties)\") into one line\n// joined by \" / \".\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the two \"Heat map\" list-item paragraphs by their exact text.\nlet statePara = null;   // \"Heat map of voting for US (by state)\"\nlet countyPara = null;  // \"Heat map for GA (by counties)\"\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (text === \"Heat map of voting for US (by state)\") {\n    statePara = para;\n  } else if (text === \"Heat map for GA (by counties)\") {\n    countyPara = para;\n  }\n}\n\nif (!statePara || !countyPara) {\n  throw new Error(\"Could not locate the expected 'Heat map' paragraphs.\");\n}\n\n// 1) Remove the stray _GoBack bookmark from wherever it currently sits.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2) Append \" / Heat map for GA (by counties)\" onto the \"(by state)\" line.\nstatePara.insertText(\" / \", Word.InsertLocation.end);\nstatePara.insertText(\"Heat map for GA (by counties)\", Word.InsertLocation.end);\n\n// 3) Empty out the now-redundant \"Heat map for GA (by counties)\" paragraph,\n//    keeping the paragraph (and its list formatting) itself intact.\nconst countyContent = countyPara.getRange(\"Content\");\ncountyContent.insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 4) Re-create the _GoBack bookmark inside that now-empty paragraph.\nconst countyEnd = countyPara.getRange(Word.RangeLocation.end);\ncountyEnd.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Move the \"_GoBack\" bookmark from the empty list item (after the\n# \"Charts we wanted to build using D3.js?\" line) down onto the empty\n# list item left behind after merging the two \"Heat map\" bullets\n# (\"... (by state)\" and \"Heat map for GA (by counties)\") into one line\n# joined by \" / \".\n\n$d = $word.ActiveDocument\n\n# Locate the two \"Heat map\" list-item paragraphs by their exact text.\n$stateParaIndex = 0\n$countyParaIndex = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"Heat map of voting for US (by state)\") {\n        $stateParaIndex = $i\n    } elseif ($t -eq \"Heat map for GA (by counties)\") {\n        $countyParaIndex = $i\n    }\n}\n\nif ($stateParaIndex -eq 0 -or $countyParaIndex -eq 0) {\n    throw \"Could not locate the expected 'Heat map' paragraphs.\"\n}\n\n# 1) Remove the stray _GoBack bookmark from wherever it currently sits.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Append \" / Heat map for GA (by counties)\" onto the \"(by state)\" line\n#    as two separate runs (mirrors the original authoring).\n$stateRange = $d.Paragraphs($stateParaIndex).Range\n$stateRange.End = $stateRange.End - 1\n$stateRange.InsertAfter(\" / \")\n$stateRange.InsertAfter(\"Heat map for GA (by counties)\")\n\n# 3) Empty out the now-redundant \"Heat map for GA (by counties)\" paragraph,\n#    keeping the paragraph (and its list formatting) itself intact.\n$countyRange = $d.Paragraphs($countyParaIndex).Range\n$countyRange.End = $countyRange.End - 1\n$countyRange.Text = \"\"\n\n# 4) Re-create the _GoBack bookmark inside that now-empty paragraph.\n$d.Bookmarks.Add(\"_GoBack\", $countyRange)\n"}
